$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.136.86"
$ws.Range("E2").Value = "  -4.48%  "
$ws.Range("D3").Value = "2.223.72"
$ws.Range("E3").Value = "  -5.75%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.02"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.93"
$ws.Range("E6").Value = "  -8.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.581"
$ws.Range("E7").Value = "  -7.36%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.563"
$ws.Range("E9").Value = "  -8.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.79"
$ws.Range("E10").Value = "  -10.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.33"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0827"
$ws.Range("E12").Value = "  -9.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.76"
$ws.Range("E13").Value = "  -8.18%  "
$ws.Range("E14").Value = "  -3.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.870"
$ws.Range("E15").Value = "  -11.31%  "
$ws.Range("D16").Value = "2.560.83"
$ws.Range("E16").Value = "  -6.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.12"
$ws.Range("E17").Value = "  -7.31%  "
$ws.Range("D18").Value = "2.224.79"
$ws.Range("E18").Value = "  -5.59%  "
$ws.Range("D19").Value = "42.983.32"
$ws.Range("E19").Value = "  -4.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.77"
$ws.Range("E20").Value = "  +4.06%  "
$ws.Range("D21").Value = "0.0₃0962"
$ws.Range("E21").Value = "  -9.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.46"
$ws.Range("E22").Value = "  -11.78%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.21"
$ws.Range("E23").Value = "  -7.99%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.11"
$ws.Range("E24").Value = "  -10.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "236.98"
$ws.Range("E25").Value = "  -8.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.13"
$ws.Range("E26").Value = "  -8.70%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  -9.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -5.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.29"
$ws.Range("E30").Value = "  -13.93%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0886"
$ws.Range("E31").Value = "  -8.07%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.52"
$ws.Range("E32").Value = "  -8.00%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.99"
$ws.Range("E33").Value = "  -10.01%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.79"
$ws.Range("E34").Value = "  -7.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.77"
$ws.Range("E35").Value = "  -5.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.28"
$ws.Range("E36").Value = "  +10.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.05"
$ws.Range("E37").Value = "  +18.57%  "
$ws.Range("E38").Value = "  -5.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.45"
$ws.Range("E39").Value = "  -7.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  -11.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.75"
$ws.Range("E41").Value = "  -4.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0325"
$ws.Range("E42").Value = "  -8.46%  "
$ws.Range("D43").Value = "1.884.41"
$ws.Range("E43").Value = "  +13.02%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.25"
$ws.Range("E45").Value = "  -4.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.01"
$ws.Range("E46").Value = "  -10.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.208"
$ws.Range("E47").Value = "  -10.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.48"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.44"
$ws.Range("E49").Value = "  -4.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.89"
$ws.Range("E50").Value = "  -12.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.75"
$ws.Range("E51").Value = "  -5.50%  "
